# Atualizado por script em 02-12-2023 14:46
#
# 1) Rows 46/47 and rows 58/59 were re-sorted: the match data held in
#    columns F:V swaps between the two rows in each pair (columns A:E -
#    index / country / tournament / season / match-date - stay put).
# 2) A brand-new match (row 61) is appended after the previous last row
#    (row 60), extending the sheet from A1:V60 to A1:V61.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param([int]$rowA, [int]$rowB)

    # Columns F..V (6..22) hold the match-specific payload that needs to
    # trade places; columns A..E are left untouched.
    for ($col = 6; $col -le 22; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)
        $valA = $cellA.Value()
        $valB = $cellB.Value()
        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}

Swap-RowData 46 47
Swap-RowData 58 59

# New row 61, appended after the existing last row (60). Copy row 60 first
# so the A61/E61 formatting (bold+border index column, datetime column)
# matches the rest of the table, then overwrite every cell with the new
# match's real values.
$ws.Range("A60:V60").Copy($ws.Range("A61:V61"))

$ws.Range("A61").Value = 60
$ws.Range("B61").Value = "united-arab-emirates"
$ws.Range("C61").Value = "uae-league"
$ws.Range("D61").Value = "2023-2024"
$ws.Range("E61").Value = 45262.57291666666
$ws.Range("F61").Value = "Al Wasl"
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = "Al Bataeh"
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 1.25
$ws.Range("K61").Value = "25/11/2023 18:13"
$ws.Range("L61").Value = 1.35
$ws.Range("M61").Value = "02/12/2023 13:36"
$ws.Range("N61").Value = 6.12
$ws.Range("O61").Value = "25/11/2023 18:13"
$ws.Range("P61").Value = 5.74
$ws.Range("Q61").Value = "02/12/2023 13:41"
$ws.Range("R61").Value = 8.33
$ws.Range("S61").Value = "25/11/2023 18:13"
$ws.Range("T61").Value = 7.37
$ws.Range("U61").Value = "02/12/2023 13:41"
$ws.Range("V61").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-wasl-al-bataeh/hdxkY2PF/"
